$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Customer Class section (rows 3-6) ---
# Row 3: Author notation -> full 1 point earned
$ws.Range("E3").Value = 1

# Row 4: Constructor -> full 2 points earned
$ws.Range("E4").Value = 2

# Row 5: Getter method -> full 2 points earned
$ws.Range("E5").Value = 2

# Row 6: toString() method -> 1 of 2 points earned, with grading comment
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = "(-1) for incorrect output in toString method"

# --- Product Class section (rows 10-14) ---
# Row 10: Two argument constructor -> full 2 points earned
$ws.Range("E10").Value = 2

# Row 11: Getter methods -> full 2 points earned
$ws.Range("E11").Value = 2

# Row 12: hashcode() method -> full 2 points earned
$ws.Range("E12").Value = 2

# Row 13: equals() method -> full 2 points earned
$ws.Range("E13").Value = 2

# Row 14: toString() method -> 1 of 2 points earned, with grading comment
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = "(-1) for incorrect output in toString method"

# Recalculate the workbook so totals (SUM formulas) refresh
$excel.CalculateFullRebuild()

# Restore the view: scroll back to top-left and set the final selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F15").Select()
